$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-18 07:55:23"
$wsZh.Range("G2").Value = "2016-02-18 07:56:08"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-18 07:55:34"
$wsDe.Range("G2").Value = "2016-02-18 07:56:28"
